# Atualizado por script em 12-11-2023 20:45
#
# The upstream scraper re-ran and the ordering of a handful of same-day
# matches changed (and one brand new match was scraped). For each affected
# pair of rows the row index (col A) and match date (col E) stay put, but
# the match details (columns F:V - teams, scores, odds, timestamps and
# match url) swap between the two rows. A new row (117) is appended for a
# match that wasn't present before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchDetails($rowA, $rowB) {
    $rangeA = $ws.Range("F$($rowA):V$($rowA)")
    $rangeB = $ws.Range("F$($rowB):V$($rowB)")
    $tmp = $rangeA.Value2
    $rangeA.Value2 = $rangeB.Value2
    $rangeB.Value2 = $tmp
}

# Row pairs whose match details (F:V) were swapped
Swap-MatchDetails 88 89
Swap-MatchDetails 97 98
Swap-MatchDetails 108 110
Swap-MatchDetails 112 113
Swap-MatchDetails 115 116

# New row 117 (copy formatting of the row above it, then fill in values)
$ws.Range("A116:V116").Copy($ws.Range("A117:V117"))

$ws.Range("A117").Value2 = 116
$ws.Range("B117").Value2 = "serbia"
$ws.Range("C117").Value2 = "prva-liga"
$ws.Range("D117").Value2 = "2023-2024"
$ws.Range("E117").Value2 = 45242.70833333334
$ws.Range("F117").Value2 = "Kolubara"
$ws.Range("G117").Value2 = 4
$ws.Range("H117").Value2 = "Vrsac"
$ws.Range("I117").Value2 = 0
$ws.Range("J117").Value2 = 1.8
$ws.Range("K117").Value2 = "12/11/2023 06:13"
$ws.Range("L117").Value2 = 1.97
$ws.Range("M117").Value2 = "12/11/2023 16:47"
$ws.Range("N117").Value2 = 3.03
$ws.Range("O117").Value2 = "12/11/2023 06:13"
$ws.Range("P117").Value2 = 2.92
$ws.Range("Q117").Value2 = "12/11/2023 16:47"
$ws.Range("R117").Value2 = 4.41
$ws.Range("S117").Value2 = "12/11/2023 06:13"
$ws.Range("T117").Value2 = 3.92
$ws.Range("U117").Value2 = "12/11/2023 16:47"
$ws.Range("V117").Value2 = "https://www.betexplorer.com/football/serbia/prva-liga/kolubara-vrsac/8WMVRxXb/"
